# "modified config of skill" - update the AtkDis (column H) values for all
# skill rows on Sheet1 from 3 to 2.5, clearing the border styling some of
# those cells had, and leave the selection on H9 (matching the final
# on-screen state after the edit).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("H2:H9")

# Clear the border formatting some of these cells carried (H4:H9 had a
# top/bottom border applied) so every AtkDis cell ends up with the plain,
# unstyled look.
$rng.Borders.LineStyle = -4142

# New AtkDis value for every skill.
$rng.Value = 2.5

# Leave the active selection on H9, matching the saved view state.
$ws.Range("H9").Select() | Out-Null
